$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(943).Insert()

$ws.Cells.Item(943, 1).Value = 3
$ws.Cells.Item(943, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(943, 3).Value = "Coquimbo"
$ws.Cells.Item(943, 4).Value = 45223
$ws.Cells.Item(943, 5).Value = 5
$ws.Cells.Item(943, 6).Value = 100112045
$ws.Cells.Item(943, 7).Value = "Zapallo"
$ws.Cells.Item(943, 8).Value = "Camote"
$ws.Cells.Item(943, 9).Value = "1a (guarda)"
$ws.Cells.Item(943, 10).Value = 100
$ws.Cells.Item(943, 11).Value = 1000
$ws.Cells.Item(943, 12).Value = 1000
$ws.Cells.Item(943, 13).Value = 1000
$ws.Cells.Item(943, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(943, 15).Value = "Provincia de Talca"
$ws.Cells.Item(943, 16).Value = 1000
$ws.Cells.Item(943, 17).Value = 1
$ws.Cells.Item(943, 18).Value = "Hortaliza"
